$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '61.782.62'
$ws.Range('E2').Value = '  +0.96%  '

# Row 3
$ws.Range('D3').Value = '3.388.65'
$ws.Range('E3').Value = '  +0.43%  '

# Row 4
$ws.Range('E4').Value = '  -0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.12%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.77'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.71%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.03%  '

# Row 8
$ws.Range('D8').Value = '3.388.16'
$ws.Range('E8').Value = '  +0.39%  '

# Row 9
$ws.Range('E9').Value = '  -0.60%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.50'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.00%  '

# Row 11
$ws.Range('E11').Value = '  +3.16%  '

# Row 12
$ws.Range('E12').Value = '  +1.31%  '

# Row 13
$ws.Range('D13').Value = '3.969.00'
$ws.Range('E13').Value = '  +0.56%  '

# Row 14
$ws.Range('E14').Value = '  +1.75%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000178'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.14%  '

# Row 16
$ws.Range('D16').Value = '3.408.25'
$ws.Range('E16').Value = '  +0.91%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.44'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.30%  '

# Row 18
$ws.Range('D18').Value = '61.874.59'
$ws.Range('E18').Value = '  +0.71%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.17'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.03%  '

# Row 20
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.83'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.69%  '

# Row 21
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.49'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.46%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '386.07'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.35%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.565'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.43%  '

# Row 24
$ws.Range('B24').Value = 'WrappedeETH'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D24').Value = '3.535.29'
$ws.Range('E24').Value = '  +0.70%  '

# Row 25
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000128'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +8.81%  '

# Row 26
$ws.Range('E26').Value = '  -0.18%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '71.46'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.05%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.69'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.51%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.63'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.63%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.04%  '

# Row 31
$ws.Range('E31').Value = '  +2.86%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.23'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.60%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.19'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.94%  '

# Row 34
$ws.Range('E34').Value = '  +0.01%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.50'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.44%  '

# Row 36
$ws.Range('D36').Value = '3.419.42'
$ws.Range('E36').Value = '  +0.42%  '

# Row 37
$ws.Range('E37').Value = '  -3.24%  '

# Row 38
$ws.Range('E38').Value = '  +2.57%  '

# Row 39
$ws.Range('E39').Value = '  -1.07%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '165.29'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.96%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0788'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.47%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.75'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.83%  '

# Row 43
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.787'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.63%  '

# Row 44
$ws.Range('B44').Value = 'ONDO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.24'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.62%  '

# Row 45
$ws.Range('E45').Value = '  -0.01%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.08'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.77%  '

# Row 47
$ws.Range('E47').Value = '  +0.45%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '41.38'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.46%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.89'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.36%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.96'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.05%  '

# Row 51
$ws.Range('D51').Value = '2.347.87'
$ws.Range('E51').Value = '  +6.97%  '
